# Weekly price update: insert two new daily-price rows at the top of the
# data block (rows 435-436) for "Vega Monumental Concepción - Papa",
# pushing the existing rows 435:450 down to 437:452.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 435; everything below shifts
# down by two (435->437 ... 450->452).
$ws.Rows("435:436").Insert()

# New row 435 — Asterix, 1a (cosecha), Región de Los Lagos
$ws.Cells.Item(435, 1).Value = 11
$ws.Cells.Item(435, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(435, 3).Value = "Bíobío"
$ws.Cells.Item(435, 4).Value = 45041
$ws.Cells.Item(435, 5).Value = 8
$ws.Cells.Item(435, 6).Value = 100114001
$ws.Cells.Item(435, 7).Value = "Papa"
$ws.Cells.Item(435, 8).Value = "Asterix"
$ws.Cells.Item(435, 9).Value = "1a (cosecha)"
$ws.Cells.Item(435, 10).Value = 220
$ws.Cells.Item(435, 11).Value = 12000
$ws.Cells.Item(435, 12).Value = 13000
$ws.Cells.Item(435, 13).Value = 12455
$ws.Cells.Item(435, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(435, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(435, 16).Value = 498
$ws.Cells.Item(435, 17).Value = 25
$ws.Cells.Item(435, 18).Value = "Hortaliza"

# New row 436 — Patagonia, 1a (cosecha), Región de La Araucanía
$ws.Cells.Item(436, 1).Value = 11
$ws.Cells.Item(436, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(436, 3).Value = "Bíobío"
$ws.Cells.Item(436, 4).Value = 45041
$ws.Cells.Item(436, 5).Value = 8
$ws.Cells.Item(436, 6).Value = 100114001
$ws.Cells.Item(436, 7).Value = "Papa"
$ws.Cells.Item(436, 8).Value = "Patagonia"
$ws.Cells.Item(436, 9).Value = "1a (cosecha)"
$ws.Cells.Item(436, 10).Value = 220
$ws.Cells.Item(436, 11).Value = 10000
$ws.Cells.Item(436, 12).Value = 11000
$ws.Cells.Item(436, 13).Value = 10545
$ws.Cells.Item(436, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(436, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(436, 16).Value = 422
$ws.Cells.Item(436, 17).Value = 25
$ws.Cells.Item(436, 18).Value = "Hortaliza"
